{"js": "// Replace the 25 two-digit multiplication prompts in the practice table\n// with the new values from the target revision. Each old prompt string is\n// unique in the document, so a simple body-wide search/replace per pair is\n// sufficient and robust to table/row layout.\nconst replacements = [\n  [\"57\u00d799=\", \"74\u00d796=\"],\n  [\"25\u00d749=\", \"61\u00d788=\"],\n  [\"14\u00d749=\", \"26\u00d747=\"],\n  [\"43\u00d789=\", \"92\u00d796=\"],\n  [\"43\u00d761=\", \"43\u00d777=\"],\n  [\"37\u00d724=\", \"98\u00d713=\"],\n  [\"36\u00d766=\", \"69\u00d723=\"],\n  [\"13\u00d771=\", \"44\u00d780=\"],\n  [\"84\u00d729=\", \"30\u00d780=\"],\n  [\"75\u00d712=\", \"17\u00d771=\"],\n  [\"31\u00d741=\", \"58\u00d771=\"],\n  [\"84\u00d764=\", \"91\u00d769=\"],\n  [\"52\u00d759=\", \"62\u00d751=\"],\n  [\"95\u00d760=\", \"15\u00d790=\"],\n  [\"48\u00d762=\", \"42\u00d754=\"],\n  [\"36\u00d746=\", \"71\u00d784=\"],\n  [\"71\u00d777=\", \"86\u00d784=\"],\n  [\"12\u00d721=\", \"59\u00d780=\"],\n  [\"47\u00d742=\", \"87\u00d728=\"],\n  [\"58\u00d785=\", \"47\u00d745=\"],\n  [\"87\u00d714=\", \"21\u00d730=\"],\n  [\"44\u00d784=\", \"84\u00d795=\"],\n  [\"96\u00d740=\", \"14\u00d766=\"],\n  [\"30\u00d714=\", \"55\u00d721=\"],\n  [\"79\u00d718=\", \"80\u00d717=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit multiplication prompts in the practice table\n# with the new values from the target revision. Each old prompt string is\n# unique in the document, so Find/Replace-All per pair is sufficient and\n# robust regardless of table/row layout.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"57\u00d799=\", \"74\u00d796=\"),\n    @(\"25\u00d749=\", \"61\u00d788=\"),\n    @(\"14\u00d749=\", \"26\u00d747=\"),\n    @(\"43\u00d789=\", \"92\u00d796=\"),\n    @(\"43\u00d761=\", \"43\u00d777=\"),\n    @(\"37\u00d724=\", \"98\u00d713=\"),\n    @(\"36\u00d766=\", \"69\u00d723=\"),\n    @(\"13\u00d771=\", \"44\u00d780=\"),\n    @(\"84\u00d729=\", \"30\u00d780=\"),\n    @(\"75\u00d712=\", \"17\u00d771=\"),\n    @(\"31\u00d741=\", \"58\u00d771=\"),\n    @(\"84\u00d764=\", \"91\u00d769=\"),\n    @(\"52\u00d759=\", \"62\u00d751=\"),\n    @(\"95\u00d760=\", \"15\u00d790=\"),\n    @(\"48\u00d762=\", \"42\u00d754=\"),\n    @(\"36\u00d746=\", \"71\u00d784=\"),\n    @(\"71\u00d777=\", \"86\u00d784=\"),\n    @(\"12\u00d721=\", \"59\u00d780=\"),\n    @(\"47\u00d742=\", \"87\u00d728=\"),\n    @(\"58\u00d785=\", \"47\u00d745=\"),\n    @(\"87\u00d714=\", \"21\u00d730=\"),\n    @(\"44\u00d784=\", \"84\u00d795=\"),\n    @(\"96\u00d740=\", \"14\u00d766=\"),\n    @(\"30\u00d714=\", \"55\u00d721=\"),\n    @(\"79\u00d718=\", \"80\u00d717=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
